# "Absenzenlisten-Templates 2016/2017 minimal ueberarbeitet
#  (zu grosse Schrift in einzelnen Zellen korrigiert)"
#
# In the student attendance table, each data row has a pair of cells
# (an empty cell immediately followed by a cell containing an "X")
# whose paragraph/run formatting is missing the w:sz/w:szCs run
# properties that every other cell in the table already carries
# (10pt / half-point value 20). That makes the "X" render oversized
# compared to the rest of the table. Fix those two cells per row by
# explicitly setting the character size to 10pt, which brings the
# run and paragraph-mark formatting in line with the surrounding
# cells (adds <w:sz w:val="20"/><w:szCs w:val="20"/> to their w:rPr).

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $row = $tbl.Rows.Item($r)
    $n = $row.Cells.Count

    for ($c = 1; $c -le $n; $c++) {
        $cell = $row.Cells.Item($c)
        $txt = $cell.Range.Text.Trim([char]13, [char]7)

        if ($txt -eq "X" -and $cell.Range.Font.Size -ne 10) {
            # The "X" marker cell itself.
            $cell.Range.Font.Size = 10
            $cell.Range.Font.SizeBi = 10

            # The empty cell immediately preceding it in the same row.
            if ($c -gt 1) {
                $prevCell = $row.Cells.Item($c - 1)
                $prevTxt = $prevCell.Range.Text.Trim([char]13, [char]7)
                if ($prevTxt -eq "") {
                    $prevCell.Range.Font.Size = 10
                    $prevCell.Range.Font.SizeBi = 10
                }
            }
        }
    }
}

Write-Output "done"
